# Mark a set of already-completed roadmap items as struck-through
# (completed-task styling), matching the "deuxieme remise" update.
#
# The document is a checklist ("Roadmap") where finished bullet items are
# shown with strikethrough formatting. This change extends the
# strikethrough formatting to the remaining bullet points that describe
# tasks which are now considered done.

$d = $word.ActiveDocument

# Exact paragraph texts that must receive strikethrough formatting.
# (Several of these texts repeat more than once in the document, so we
# match by paragraph content rather than by a single global search.)
$targets = @(
    "Le montant total des revenus",
    "Le montant total des dépenses",
    "Faire le calcul en conséquence de la fréquence (exemple si bimensuel, faire fois 2)",
    "Sauvegarder le montant quelque part (dans une liste peut-être)",
    "Créer l’action de sauvegarder les informations et de retourner à l’activité principale quand le bouton d’ajout est cliqué",
    "Additionner le revenu avec le total des revenus déjà existants",
    "Additionner la dépense avec le total des dépenses déjà existantes",
    "Afficher le montant total des revenus additionnés",
    "Afficher le montant total des dépenses additionnées",
    "Faire le calcul des revenus moins les dépenses et sauvegarder le montant",
    "Créer un espace au milieu/bas contenant",
    "Un label d’argent disponible",
    "Le montant total d’argent disponible"
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.TrimEnd("`r", "`n", "`a").Trim()

    foreach ($target in $targets) {
        if ($text -eq $target -or $text.StartsWith($target)) {
            $p.Range.Font.StrikeThrough = $true
            break
        }
    }
}

Write-Host "Strikethrough formatting applied."
